$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows with the new questionnaire text for the
# Cost Application portal screen.
$ws.Range("A22").Value = "Cost response received"
$ws.Range("A23").Value = "Cost response reference number"

# Keep selection/active cell consistent with where Excel would leave it
# after typing into the new last row.
$ws.Range("A23").Select()
